# New changes in ops (ordercreation & orderpage & order form)
# 12/30/2024 - insert "Typist" / "Typist QC" columns and refresh order data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert two new columns at E:F (old E..M shift right to G..O) and
#    give them their new headers + values.
# ---------------------------------------------------------------------
$ws.Range("E1:F1").EntireColumn.Insert()

$ws.Range("E1").Value = "Typist"
$ws.Range("F1").Value = "Typist QC"

$ws.Range("E2").Value = "SIPL5317"
$ws.Range("F2").Value = "SIPL5317"

$ws.Range("E3").Value = "SIPL0102"
$ws.Range("F3").Value = "SIPL0103"

# ---------------------------------------------------------------------
# 2. Data corrections on the shifted-right columns.
# ---------------------------------------------------------------------
# Row 2: the trailing "Search(T1)" tag is cleared out.
$ws.Range("O2").Value = ""

# Row 3: Emp ID / Assignee_QA no longer carry values.
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""

# Row 3: Status moves from "WIP" to "Typing", and the tail tag becomes
# "Typing(T1)" instead of "Search(T1)".
$ws.Range("N3").Value = "Typing"
$ws.Range("O3").Value = "Typing(T1)"

# ---------------------------------------------------------------------
# 3. Column width touch-ups to match the refreshed layout.
# ---------------------------------------------------------------------
$ws.Range("E1").ColumnWidth = 7.5
$ws.Range("H1").ColumnWidth = 10.166666666666666
$ws.Range("J1").ColumnWidth = 15.666666666666666

# ---------------------------------------------------------------------
# 4. Leave the cursor where the author ended up.
# ---------------------------------------------------------------------
$ws.Range("J7").Select()
